$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 1200
$ws.Range("J17").Value = 1200
$ws.Range("L17").Value = 3600
$ws.Range("N17").Value = -3936

# Row 58
$ws.Range("H58").Value = 786.25
$ws.Range("I58").Value = 73.75
$ws.Range("J58").Value = 1498.75
$ws.Range("K58").Value = 221.25
$ws.Range("L58").Value = 4496.25
$ws.Range("M58").Value = -71.25
$ws.Range("N58").Value = -4796.25

# Row 80
$ws.Range("H80").Value = 2041.5834
$ws.Range("J80").Value = 3999.75
$ws.Range("L80").Value = 11999.25
$ws.Range("N80").Value = -13995.25

# Row 83
$ws.Range("H83").Value = 2041.5834
$ws.Range("J83").Value = 3999.75
$ws.Range("L83").Value = 35997.75
$ws.Range("N83").Value = -45981.75

# Row 135
$ws.Range("H135").Value = 2342.3333
$ws.Range("I135").Value = 2342.3333
$ws.Range("K135").Value = 21080.9997
$ws.Range("M135").Value = -18545.9997

# ---------------------------------------------------------------------------
# Sheet ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 8
$ws.Range("H8").Value = 1433375
$ws.Range("I8").Value = 1542711.5
$ws.Range("K8").Value = 1542711.5
$ws.Range("M8").Value = -1542567.5

# Row 31
$ws.Range("H31").Value = 15500
$ws.Range("I31").Value = 15500
$ws.Range("K31").Value = 15500
$ws.Range("M31").Value = -15206

# Row 74
$ws.Range("H74").Value = 946.5714
$ws.Range("I74").Value = 990.5
$ws.Range("J74").Value = 888
$ws.Range("K74").Value = 990.5
$ws.Range("L74").Value = 888
$ws.Range("M74").Value = -116.5
$ws.Range("N74").Value = -2636

# Row 77
$ws.Range("H77").Value = 946.5714
$ws.Range("I77").Value = 990.5
$ws.Range("J77").Value = 888
$ws.Range("K77").Value = 4952.5
$ws.Range("L77").Value = 4440
$ws.Range("M77").Value = -584.5
$ws.Range("N77").Value = -13176

# Row 97
$ws.Range("H97").Value = 6797.3335
$ws.Range("I97").Value = 147.125
$ws.Range("K97").Value = 147.125
$ws.Range("M97").Value = 348.875

# Row 102
$ws.Range("H102").Value = 2259.8
$ws.Range("I102").Value = 799.8570999999999
$ws.Range("K102").Value = 799.8570999999999
$ws.Range("M102").Value = 822.1429000000001

# Row 110
$ws.Range("H110").Value = 1403.8182
$ws.Range("I110").Value = 1344.2
$ws.Range("K110").Value = 1344.2
$ws.Range("M110").Value = 700.8

# Row 132
$ws.Range("H132").Value = 2646.6
$ws.Range("I132").Value = 2805.25
$ws.Range("K132").Value = 8415.75
$ws.Range("M132").Value = -5885.75

# ---------------------------------------------------------------------------
# Sheet BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 11
$ws.Range("H11").Value = 590.75
$ws.Range("I11").Value = 121
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 121
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 19
$ws.Range("N11").Value = -2280

# Row 20
$ws.Range("H20").Value = 3339.8333
$ws.Range("I20").Value = 2054.1428
$ws.Range("J20").Value = 5139.8
$ws.Range("K20").Value = 2054.1428
$ws.Range("L20").Value = 5139.8
$ws.Range("M20").Value = -1807.1428
$ws.Range("N20").Value = -5633.8

# Row 99
$ws.Range("H99").Value = 2091.25
$ws.Range("I99").Value = 1955.625
$ws.Range("K99").Value = 1955.625
$ws.Range("M99").Value = -457.625

# Row 102
$ws.Range("H102").Value = 12725
$ws.Range("I102").Value = 12725
$ws.Range("K102").Value = 12725
$ws.Range("M102").Value = -9480

# Row 134
$ws.Range("H134").Value = 2048.2068
$ws.Range("I134").Value = 1592.2307
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 4776.6921
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -2241.6921
$ws.Range("N134").Value = -23070

# ---------------------------------------------------------------------------
# Sheet CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 4659.857
$ws.Range("I31").Value = 2922.8
$ws.Range("K31").Value = 2922.8
$ws.Range("M31").Value = -2627.8

# Row 34
$ws.Range("H34").Value = 4659.857
$ws.Range("I34").Value = 2922.8
$ws.Range("K34").Value = 2922.8
$ws.Range("M34").Value = -2720.8

# Row 107
$ws.Range("I107").Value = 350.46155
$ws.Range("J107").Value = 999.3333
$ws.Range("K107").Value = 350.46155
$ws.Range("L107").Value = 999.3333
$ws.Range("M107").Value = 1569.53845
$ws.Range("N107").Value = -4839.3333

# ---------------------------------------------------------------------------
# Sheet CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 121
$ws.Range("H121").Value = 1838.6
$ws.Range("J121").Value = 2998
$ws.Range("L121").Value = 8994
$ws.Range("N121").Value = -11614

# Row 131
$ws.Range("H131").Value = 1092.3334
$ws.Range("I131").Value = 666.6667
$ws.Range("J131").Value = 1305.1666
$ws.Range("K131").Value = 2000.0001
$ws.Range("L131").Value = 3915.4998
$ws.Range("M131").Value = 3039.9999
$ws.Range("N131").Value = -13995.4998

# ---------------------------------------------------------------------------
# Sheet GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value = 1261.4445
$ws.Range("I97").Value = 1411.5
$ws.Range("J97").Value = 961.3333
$ws.Range("K97").Value = 1411.5
$ws.Range("L97").Value = 961.3333
$ws.Range("M97").Value = -915.5
$ws.Range("N97").Value = -1953.3333

# Row 99
$ws.Range("H99").Value = 10734
$ws.Range("I99").Value = 10734
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10734
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8488
$ws.Range("N99").Value = ""

# Row 132
$ws.Range("H132").Value = 2539.4285
$ws.Range("I132").Value = 1842.1765
$ws.Range("K132").Value = 5526.529500000001
$ws.Range("M132").Value = -2996.529500000001

# ---------------------------------------------------------------------------
# Sheet LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 398.33334
$ws.Range("I22").Value = 398.33334
$ws.Range("K22").Value = 398.33334
$ws.Range("M22").Value = -103.33334

# Row 27
$ws.Range("H27").Value = 398.33334
$ws.Range("I27").Value = 398.33334
$ws.Range("K27").Value = 398.33334
$ws.Range("M27").Value = -291.33334

# Row 46
$ws.Range("H46").Value = 1973.7222
$ws.Range("I46").Value = 1725
$ws.Range("J46").Value = 2172.7
$ws.Range("K46").Value = 1725
$ws.Range("L46").Value = 2172.7
$ws.Range("M46").Value = -1537
$ws.Range("N46").Value = -2548.7

# Row 93
$ws.Range("H93").Value = 1289.3334
$ws.Range("I93").Value = 1003
$ws.Range("J93").Value = 1371.1428
$ws.Range("K93").Value = 1003
$ws.Range("L93").Value = 1371.1428
$ws.Range("M93").Value = 245
$ws.Range("N93").Value = -3867.1428

# Row 99
$ws.Range("H99").Value = 90259
$ws.Range("I99").Value = 90259
$ws.Range("K99").Value = 90259
$ws.Range("M99").Value = -87264

# Row 136
$ws.Range("H136").Value = 2950.6
$ws.Range("I136").Value = 2562.625
$ws.Range("J136").Value = 4502.5
$ws.Range("K136").Value = 7687.875
$ws.Range("L136").Value = 13507.5
$ws.Range("M136").Value = -5137.875
$ws.Range("N136").Value = -18607.5

# ---------------------------------------------------------------------------
# Sheet WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""

# Row 132
$ws.Range("H132").Value = 101353.9
$ws.Range("I132").Value = 144391.28
$ws.Range("K132").Value = 433173.84
$ws.Range("M132").Value = -430643.84
